# Invincibility powerup logic:
# 1) Drop a "_GoBack" bookmark (the marker Word leaves at the last edit
#    point) around the "Also a timer ... Lose the Game." sentence.
# 2) The old "_GoBack" bookmark used to sit inside the "Exit Fence" bullet,
#    splitting it into two runs ("...upwards and" / " disappears."). Re-join
#    that bullet into a single run with the full sentence; replacing across
#    the split also clears out the stray bookmark that was living there.

$d = $word.ActiveDocument

# --- Change 1: bookmark the "Lose the Game" sentence -----------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute(
    "Also a timer where the Gas Station Fuel pump explodes and you Lose the Game.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $d.Bookmarks.Add("_GoBack", $rng1)
}

# --- Change 2: re-merge the "Exit Fence" bullet into a single run ----------
$rng2 = $d.Content
$rng2.Find.Execute(
    "Fence Translates upwards and disappears.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Fence Translates upwards and disappears.", 2)
